# From v1.1.1 to v1.2
# The "TC2" test case previously described "cancelar diária" and the
# "TC3" test case previously described "analisar prestação de contas".
# This edit swaps the Steps / Expected Results content between the two
# test cases so TC2 now covers "analisar prestação de contas" and TC3
# now covers "cancelar diária" (all other labels/rows stay unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC2 steps row (row 18): B18 = Steps, D18 = Expected Results
$tc2Steps  = $ws.Range("B18").Value2
$tc2Result = $ws.Range("D18").Value2

# TC3 steps row (row 25): B25 = Steps, D25 = Expected Results
$tc3Steps  = $ws.Range("B25").Value2
$tc3Result = $ws.Range("D25").Value2

# Swap them
$ws.Range("B18").Value2 = $tc3Steps
$ws.Range("D18").Value2 = $tc3Result

$ws.Range("B25").Value2 = $tc2Steps
$ws.Range("D25").Value2 = $tc2Result
